# BBI-23-1.xlsx update
# - Adjusts several students' grade inputs (columns B, E, G) on the "Лист1" sheet.
# - Two students' КР-2 column (E) switch from a numeric score to a text remark
#   ("5(без допуска)" / "4 (без допуска)"), which turns their dependent
#   MAX/IF formulas (J/K/L) into #VALUE! errors, exactly as a real grade sheet
#   would behave once a numeric cell is replaced by explanatory text.
# - The scratch/code-snippet rows below the table (30, 32, 33 containing
#   "arr", "i", "j") are removed, shrinking the used range back to A1:S25.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update individual grade cells -----------------------------------------

# Row 3
$ws.Range("E3").Value = 3

# Row 6 - КР-2 column becomes a text remark instead of a numeric score
$ws.Range("E6").Value = "5(без допуска)"

# Row 11
$ws.Range("E11").Value = -2

# Row 13
$ws.Range("E13").Value = -2

# Row 14
$ws.Range("B14").Value = 3
$ws.Range("E14").Value = 4

# Row 17
$ws.Range("E17").Value = 4
$ws.Range("G17").Value = 5

# Row 18
$ws.Range("G18").Value = 5

# Row 21 - КР-2 column becomes a text remark instead of a numeric score
$ws.Range("E21").Value = "4 (без допуска)"
$ws.Range("G21").Value = 4

# Row 23
$ws.Range("E23").Value = 4

# Row 24
$ws.Range("G24").Value = 5

# --- Widen the new КР-2 (column E) now that it holds longer text values ----
$ws.Columns("E:E").ColumnWidth = 13.14

# --- Remove the leftover scratch rows (arr / i / j) below the table --------
$ws.Rows("30:33").Delete()

# --- Restore the selected cell as left by the author ------------------------
$ws.Range("E12").Select()
